$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first blank templated data row (row 3), shifting rows 4-31 up to 3-30.
$ws.Range("A3:E3").EntireRow.Delete()

# Remove the now-superfluous blank templated rows, keeping only the new row 3.
$ws.Range("A4:E30").EntireRow.Delete()

# Clear the thin-border template formatting from row 3 so it goes back to plain/default style.
$ws.Range("A3:E3").ClearFormats()

# Fill in the actual data row.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "martha"
$ws.Range("C3").Value = "martha"
$ws.Range("D3").Value = "ISC"
$ws.Range("E3").Value = 20122423

Write-Host "Data entered"

$null = $ws.Range("E3").Select()

# Shrink the print area to match the new, smaller data range.
$ws.PageSetup.PrintArea = '$A$1:$E$3'

# The 13 small placeholder pictures were anchored starting at (the now-deleted)
# row 4; shift them up one row so they re-anchor against row 3 like the rest
# of the sheet did when that row was removed.
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
  $shp = $ws.Shapes.Item($i)
  if ($shp.Top -gt 100 -and $shp.Top -lt 110) {
    $shp.Top = 89.25
  }
}

